# Adds 114 new scraped-error rows to Sheet1 (url in col A, msg in col B),
# continuing after the existing last data row (2032), matching the row-number
# gaps present in the source data (blank separator rows are skipped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shortMsg = "UiPath.UIAutomationNext.Exceptions.NodeNotFoundException: Could not find the user-interface (UI) element for this action.`n`nPossible solutions:`n •  Ensure application is opened and the UI element is visible on the screen at execution time`n •  Edit the Target of the UI activity and use Validation to debug the issue.`n •  If needed, re-indicate the element as its properties might have changed`n •  Use `"Check state`" activity to check the application state before executing the action`n •  Increase the `"Delay before`" value to allow time to the application to render entirely and become responsive`n   at UiPath.UIAutomationNext.Activities.TargetCommonLogic.GetSearchResultAsync(IRuntimeContext runtimeContext, ITargetAnchorable target, CancellationToken token)`n   at UiPath.UIAutomationNext.Activities.TargetBase.SearchNodeAsync(IRuntimeContext runtimeContext, ITargetAnchorable target, CancellationToken token)`n   at UiPath.UIAutomationNext.Activities.NGetText.ExecuteAsync(AsyncCodeActivityContext context, CancellationToken token)`n   at UiPath.Shared.Activities.AsyncTaskCodeActivityImplementation.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at UiPath.Shared.Activities.AsyncTaskCodeActivity.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at UiPath.Shared.Activities.ContinuableAsyncCodeActivity.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at System.Activities.AsyncCodeActivity.CompleteAsyncCodeActivityData.CompleteAsyncCodeActivityWorkItem.Execute(ActivityExecutor executor, BookmarkManager bookmarkManager)"
$longMsg = "UiPath.UIAutomationNext.Exceptions.NodeNotFoundException: Could not find the user-interface (UI) element for this action.`n`nPossible solutions:`n •  Ensure application is opened and the UI element is visible on the screen at execution time`n •  Edit the Target of the UI activity and use Validation to debug the issue.`n •  If needed, re-indicate the element as its properties might have changed`n •  Use `"Check state`" activity to check the application state before executing the action`n •  Increase the `"Delay before`" value to allow time to the application to render entirely and become responsive`n   at UiPath.UIAutomationNext.Activities.TargetCommonLogic.GetSearchResultAsync(IRuntimeContext runtimeContext, ITargetAnchorable target, CancellationToken token)`n   at UiPath.UIAutomationNext.Activities.TargetBase.SearchNodeAsync(IRuntimeContext runtimeContext, ITargetAnchorable target, CancellationToken token)`n   at UiPath.UIAutomationNext.Activities.NGetText.ExecuteAsync(AsyncCodeActivityContext context, CancellationToken token)`n   at UiPath.Shared.Activities.AsyncTaskCodeActivityImplementation.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at UiPath.Shared.Activities.AsyncTaskCodeActivity.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at UiPath.Shared.Activities.ContinuableAsyncCodeActivity.EndExecute(AsyncCodeActivityContext context, IAsyncResult result)`n   at System.Activities.AsyncCodeActivity.System.Activities.IAsyncCodeActivity.FinishExecution(AsyncCodeActivityContext context, IAsyncResult result)`n   at System.Activities.AsyncCodeActivity.CompleteAsyncCodeActivityData.CompleteAsyncCodeActivityWorkItem.Execute(ActivityExecutor executor, BookmarkManager bookmarkManager)"

$rows = @(
    @(2034, "https://www.therealreal.com/products/women/shoes/sandals/chanel-vintage-tweed-pattern-slingback-sandals-fzwmv", "short")
    ,@(2035, "https://www.therealreal.com/products/women/handbags/totes/chanel-caviar-medallion-tote-fym5f", "short")
    ,@(2036, "https://www.therealreal.com/products/women/handbags/totes/louis-vuitton-antigua-cabas-gm-fzuif", "short")
    ,@(2037, "https://www.therealreal.com/products/women/clothing/pants/yves-saint-laurent-rive-gauche-vintage-wide-leg-pants-fxu1h", "short")
    ,@(2038, "https://www.therealreal.com/products/women/handbags/handle-bags/fendi-vintage-mini-pochette-g0akv", "short")
    ,@(2039, "https://www.therealreal.com/products/women/accessories/belts/chanel-vintage-2005-waist-belt-g0mye", "short")
    ,@(2040, "https://www.therealreal.com/products/women/accessories/belts/chanel-vintage-1990-chain-link-belt-fx0q4", "short")
    ,@(2041, "https://www.therealreal.com/products/women/clothing/coats/chanel-vintage-2002-fur-coat-fy8lb", "short")
    ,@(2042, "https://www.therealreal.com/products/women/clothing/coats/valentino-vintage-2000-s-down-coat-fynfh", "short")
    ,@(2043, "https://www.therealreal.com/products/women/clothing/tops/jean-paul-gaultier-vintage-2004-tunic-g0se7", "short")
    ,@(2044, "https://www.therealreal.com/products/women/clothing/skirts/d-g-vintage-mini-skirt-frh3u", "short")
    ,@(2045, "https://www.therealreal.com/products/women/accessories/wallets/chanel-vintage-2006-2008-compact-wallet-fw9d6", "short")
    ,@(2046, "https://www.therealreal.com/products/women/clothing/coats/burberry-s-vintage-trench-coat-fywzu", "short")
    ,@(2047, "https://www.therealreal.com/products/women/handbags/handle-bags/louis-vuitton-vintage-monogram-ellipse-pm-fm1dq", "short")
    ,@(2048, "https://www.therealreal.com/products/women/clothing/skirts/gucci-vintage-midi-length-skirt-fl3ri", "short")
    ,@(2049, "https://www.therealreal.com/products/women/clothing/suits-and-sets/issey-miyake-vintage-late-1980-s-early-1990-s-skirt-set-ftndc", "short")
    ,@(2050, "https://www.therealreal.com/products/women/clothing/pants/gianni-versace-vintage-straight-leg-pants-e8qk5", "short")
    ,@(2051, "https://www.therealreal.com/products/women/clothing/jackets/dolce-gabbana-vintage-late-1990-s-early-2000-s-fur-jacket-e4xas", "short")
    ,@(2053, "https://www.therealreal.com/products/beauty/bath-and-body/body-oils/everyday-oil-mainstay-blend-8-oz-dk50r", "short")
    ,@(2054, "https://www.therealreal.com/products/beauty/bath-and-body/body-wash/joanna-vargas-vitamin-c-face-wash-enk76", "short")
    ,@(2056, "https://www.therealreal.com/products/beauty/skincare/tools-and-accessories/well-kept-exfoliating-washcloth-e97f4", "short")
    ,@(2057, "https://www.therealreal.com/products/beauty/bath-and-body/hand-soap-and-moisturizers/grown-alchemist-hand-wash-sweet-orange-cedarwood-sage-300ml-eo44i", "short")
    ,@(2058, "https://www.therealreal.com/products/beauty/bath-and-body/hand-soap-and-moisturizers/nopalera-cactus-soap-in-planta-futura-e2wc7", "short")
    ,@(2059, "https://www.therealreal.com/products/beauty/hair-care/shampoo-and-conditioner/ceremonia-guava-rescue-spray-dgdke", "short")
    ,@(2060, "https://www.therealreal.com/products/beauty/fragrance/spicy/boy-smells-cowboy-kush-eau-de-parfum-etg55", "short")
    ,@(2061, "https://www.therealreal.com/products/beauty/fragrance/spicy/boy-smells-hinoki-fantome-eau-de-parfum-etfvo", "short")
    ,@(2062, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-desire-luxury-lip-tint-elhwk", "short")
    ,@(2063, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-bare-luxury-lip-tint-elhyb", "short")
    ,@(2064, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-muse-luxury-lip-tint-elht2", "short")
    ,@(2065, "https://www.therealreal.com/products/beauty/skincare/tools-and-accessories/lanshin-rose-quartz-sculpting-spoon-eml2z", "short")
    ,@(2066, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-intrigue-luxury-lip-tint-elhv2", "short")
    ,@(2067, "https://www.therealreal.com/products/beauty/hair-care/brushes-and-combs/machete-no-4-comb-orchid-ekry2", "short")
    ,@(2068, "https://www.therealreal.com/products/beauty/skincare/ubeauty-resurfacing-compound-adrj9", "short")
    ,@(2069, "https://www.therealreal.com/products/beauty/makeup/lips/dr-devgan-platinum-lip-plump-spf-30-do7mw", "short")
    ,@(2070, "https://www.therealreal.com/products/beauty/fragrance/woody/cra-yon-sand-service-50ml-eau-de-parfum-ef3ve", "short")
    ,@(2071, "https://www.therealreal.com/products/beauty/fragrance/fresh/cra-yon-vanilla-ceo-50ml-eau-de-parfum-ef3su", "short")
    ,@(2072, "https://www.therealreal.com/products/beauty/skincare/masks-and-exfoliators/nopalera-cactus-flower-exfoliant-e2wa6", "short")
    ,@(2073, "https://www.therealreal.com/products/beauty/bath-and-body/body-moisturizers/lauren-s-all-purpose-classic-jar-e6p57", "short")
    ,@(2074, "https://www.therealreal.com/products/beauty/skincare/serums/joanna-vargas-rescue-serum-ah655", "short")
    ,@(2075, "https://www.therealreal.com/products/beauty/bath-and-body/body-wash/corpus-natural-body-wash-no-green-c9kmq", "short")
    ,@(2076, "https://www.therealreal.com/products/beauty/bath-and-body/hand-soap-and-moisturizers/nopalera-cactus-soap-in-flor-de-mayo-e2wbx", "short")
    ,@(2077, "https://www.therealreal.com/products/beauty/fragrance/spicy/19-69-chinese-tobacco-eau-de-parfum-bjpbe", "short")
    ,@(2078, "https://www.therealreal.com/products/beauty/skincare/joanna-vargas-daily-serum-9szi8", "short")
    ,@(2079, "https://www.therealreal.com/products/beauty/fragrance/woody/19-69-rainbow-bar-eau-de-parfum-bjqcn", "short")
    ,@(2080, "https://www.therealreal.com/products/beauty/hair-care/hair-masks/ceremonia-pequi-curl-activator-200ml-e4nbv", "short")
    ,@(2081, "https://www.therealreal.com/products/beauty/bath-and-body/mason-pearson-popular-mixture-hair-brush-c7ilb", "short")
    ,@(2082, "https://www.therealreal.com/products/beauty/skincare/face-oils/rowse-rosehip-oil-dt2o3", "short")
    ,@(2083, "https://www.therealreal.com/products/beauty/bath-and-body/body-oils/kindred-black-a-woman-is-fire-all-natural-perfume-dk4r8", "short")
    ,@(2084, "https://www.therealreal.com/products/beauty/bath-and-body/hand-soap-and-moisturizers/claus-porto-gift-box-9-deco-soaps-c7j30", "short")
    ,@(2085, "https://www.therealreal.com/products/beauty/skincare/masks-and-exfoliators/knc-beauty-star-eye-mask-set-dan7g", "short")
    ,@(2086, "https://www.therealreal.com/products/beauty/hair-care/crown-affair-scrunchie-no-001-pack-of-3-do7rh", "short")
    ,@(2087, "https://www.therealreal.com/products/beauty/skincare/face-oils/mara-algae-moringa-universal-face-oil-adv37", "short")
    ,@(2088, "https://www.therealreal.com/products/beauty/skincare/masks-and-exfoliators/joanna-vargas-forever-glow-anti-aging-face-mask-anmgs", "short")
    ,@(2089, "https://www.therealreal.com/products/beauty/skincare/tools-and-accessories/sounds-gua-sha-rose-quartz-doj67", "short")
    ,@(2090, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-bare-luxury-lip-tint-c9ixw", "short")
    ,@(2091, "https://www.therealreal.com/products/beauty/hair-care/shampoo-and-conditioner/ceremonia-champu-de-yucca-witch-hazel-dgdns", "short")
    ,@(2092, "https://www.therealreal.com/products/beauty/hair-care/shampoo-and-conditioner/ceremonia-guava-leave-in-conditioner-dgdfw", "short")
    ,@(2093, "https://www.therealreal.com/products/beauty/skincare/moisturizers/ubeauty-the-super-smart-hydrator-adlq7", "short")
    ,@(2094, "https://www.therealreal.com/products/beauty/makeup/face/saie-sunglow-glowy-super-gel-luminizer-cph58", "short")
    ,@(2095, "https://www.therealreal.com/products/beauty/hair-care/brushes-and-combs/crown-affair-brush-no-001-c9j4f", "short")
    ,@(2096, "https://www.therealreal.com/products/beauty/bath-and-body/body-wash/austin-austin-neroli-petitgrain-body-soap-ag1yl", "short")
    ,@(2097, "https://www.therealreal.com/products/beauty/makeup/olio-e-osso-balm-no-3-crimson-adpys", "short")
    ,@(2098, "https://www.therealreal.com/products/beauty/fragrance/ormaie-l-ivree-bleue-edp-50ml-c4deb", "short")
    ,@(2099, "https://www.therealreal.com/products/beauty/bath-and-body/body-wash/corpus-natural-body-wash-third-rose-c9k6c", "short")
    ,@(2100, "https://www.therealreal.com/products/beauty/makeup/face/saie-dreamy-liquid-blush-cph6d", "short")
    ,@(2101, "https://www.therealreal.com/products/beauty/skincare/tools-and-accessories/gilded-body-the-marble-body-brush-calacatta-viola-d8yk7", "short")
    ,@(2102, "https://www.therealreal.com/products/beauty/skincare/tools-and-accessories/gilded-body-the-marble-body-brush-lichen-flower-d8yo5", "short")
    ,@(2103, "https://www.therealreal.com/products/beauty/hair-care/brushes-and-combs/crown-affair-the-comb-no-002-d7cr7", "short")
    ,@(2104, "https://www.therealreal.com/products/beauty/hair-care/hair-masks/crown-affair-the-renewal-mask-d7cuo", "short")
    ,@(2105, "https://www.therealreal.com/products/beauty/bath-and-body/hand-soap-and-moisturizers/austin-austin-palmarosa-vetiver-hand-soap-c9kru", "short")
    ,@(2106, "https://www.therealreal.com/products/beauty/fragrance/fresh/d-s-durga-i-don-t-know-what-eau-de-parfum-100ml-ci9r9", "short")
    ,@(2107, "https://www.therealreal.com/products/beauty/fragrance/fresh/d-s-durga-i-don-t-know-what-eau-de-parfum-50ml-cia1y", "short")
    ,@(2108, "https://www.therealreal.com/products/beauty/fragrance/floral/d-s-durga-jazmin-yucatan-eau-de-parfum-50ml-ci9z2", "short")
    ,@(2109, "https://www.therealreal.com/products/beauty/makeup/face/saie-rosy-liquid-blush-cphs2", "short")
    ,@(2110, "https://www.therealreal.com/products/beauty/makeup/face/saie-the-big-makeup-brush-cpgi8", "short")
    ,@(2111, "https://www.therealreal.com/products/beauty/makeup/face/saie-starglow-glowy-super-gel-luminizer-cpife", "short")
    ,@(2112, "https://www.therealreal.com/products/beauty/makeup/lips/henne-organics-desire-luxury-lip-tint-d371r", "short")
    ,@(2113, "https://www.therealreal.com/products/beauty/fragrance/floral/regime-des-fleurs-chloe-sevigny-little-flower-eau-de-parfum-ca6e3", "short")
    ,@(2114, "https://www.therealreal.com/products/beauty/fragrance/woody/d-s-durga-debaser-pocket-perfume-cxoi5", "short")
    ,@(2115, "https://www.therealreal.com/products/beauty/bath-and-body/body-wash/flamingo-estate-body-wash-ci6bw", "short")
    ,@(2116, "https://www.therealreal.com/products/beauty/skincare/moisturizers/kindred-black-damiana-aphrodisiac-lip-skin-balm-cibes", "short")
    ,@(2117, "https://www.therealreal.com/products/beauty/hair-care/brushes-and-combs/crown-affair-the-comb-no-001-ck67i", "short")
    ,@(2118, "https://www.therealreal.com/products/beauty/hair-care/brushes-and-combs/crown-affair-the-comb-no-002-bzc9c", "short")
    ,@(2120, "https://www.therealreal.com/products/women/accessories/hats/prada-puffer-trapper-hat-elj8v", "long")
    ,@(2121, "https://www.therealreal.com/products/women/accessories/sunglasses/jeremy-scott-x-linda-farrow-shield-tinted-sunglasses-eu3m0", "long")
    ,@(2122, "https://www.therealreal.com/products/women/clothing/coats/vetements-down-coat-w-tags-bvn4p", "long")
    ,@(2123, "https://www.therealreal.com/products/women/clothing/jumpsuits-and-rompers/emilio-pucci-vintage-late-1960-s-early-1970-s-jumpsuit-esw9y", "long")
    ,@(2124, "https://www.therealreal.com/products/women/accessories/winter-accessories/prada-woven-gloves-eodj2", "short")
    ,@(2125, "https://www.therealreal.com/products/women/accessories/hats/bogner-embroidered-ski-hat-e6m69", "short")
    ,@(2126, "https://www.therealreal.com/products/women/accessories/scarves-and-shawls/loewe-paula-s-ibiza-printed-scarf-evl4a", "short")
    ,@(2127, "https://www.therealreal.com/products/women/clothing/jackets/bottega-veneta-down-jacket-fyvl7", "short")
    ,@(2128, "https://www.therealreal.com/products/women/accessories/sunglasses/celine-vintage-shield-sunglasses-foza1", "short")
    ,@(2130, "https://www.therealreal.com/products/women/shoes/boots/louis-vuitton-pillow-comfort-lv-monogram-snow-boots-d60xd", "short")
    ,@(2131, "https://www.therealreal.com/products/women/clothing/jackets/the-north-face-jacket-g1o2t", "short")
    ,@(2132, "https://www.therealreal.com/products/home/sports/snow-gear/burton-family-tree-straight-chuter-camber-splitboard-cnw30", "short")
    ,@(2133, "https://www.therealreal.com/products/women/clothing/jackets/stand-studio-down-jacket-w-tags-ez7dh", "short")
    ,@(2134, "https://www.therealreal.com/products/women/clothing/coats/chanel-vintage-1990-performance-coat-fku5t", "short")
    ,@(2135, "https://www.therealreal.com/products/women/accessories/sunglasses/cartier-oversize-tinted-sunglasses-f1kvq", "short")
    ,@(2136, "https://www.therealreal.com/products/women/accessories/winter-accessories/saint-laurent-leather-winter-gloves-g1j91", "short")
    ,@(2137, "https://www.therealreal.com/products/women/clothing/knitwear/fendi-v-neck-sweater-fusdc", "short")
    ,@(2138, "https://www.therealreal.com/products/women/accessories/hats/chanel-wool-beanie-w-tags-ff4cw", "short")
    ,@(2139, "https://www.therealreal.com/products/women/accessories/winter-accessories/celine-winter-gloves-ff9gl", "short")
    ,@(2140, "https://www.therealreal.com/products/women/handbags/luggage-and-travel/louis-vuitton-monogram-summer-trunks-keepall-bandouliere-50-djox5", "short")
    ,@(2141, "https://www.therealreal.com/products/women/clothing/jackets/emilio-pucci-vintage-2000-s-jacket-fyr4w", "short")
    ,@(2142, "https://www.therealreal.com/products/women/accessories/hats/chanel-2022-cashmere-cc-beanie-w-tags-funzo", "short")
    ,@(2143, "https://www.therealreal.com/products/women/accessories/hats/the-elder-statesman-cashmere-beanie-fxtpd", "short")
    ,@(2144, "https://www.therealreal.com/products/women/accessories/scarves-and-shawls/hermes-casaque-optique-cashmere-muffler-fw3vj", "short")
    ,@(2145, "https://www.therealreal.com/products/women/handbags/backpacks/chanel-coco-neige-shearling-backpack-erpye", "short")
    ,@(2146, "https://www.therealreal.com/products/women/clothing/tops/eckhaus-latta-floral-print-turtleneck-crop-top-g126o", "short")
    ,@(2147, "https://www.therealreal.com/products/women/accessories/scarves-and-shawls/prada-nylon-printed-scarf-fsqfz", "short")
    ,@(2148, "https://www.therealreal.com/products/women/accessories/sunglasses/celine-shield-gradient-sunglasses-g0y15", "short")
    ,@(2149, "https://www.therealreal.com/products/women/clothing/jackets/moncler-grenoble-houndstooth-print-down-jacket-fzr4v", "short")
    ,@(2150, "https://www.therealreal.com/products/home/sports/snow-gear/head-wc-rebels-racing-skis-c7brq", "short")
    ,@(2151, "https://www.therealreal.com/products/women/accessories/winter-accessories/chanel-2020-cashmere-fingerless-gloves-w-tags-ftvep", "short")
)

foreach ($row in $rows) {
    $rowNum = $row[0]
    $url = $row[1]
    $variant = $row[2]
    if ($variant -eq "long") { $msg = $longMsg } else { $msg = $shortMsg }
    $ws.Cells.Item($rowNum, 1).Value = $url
    $ws.Cells.Item($rowNum, 2).Value = $msg
}
